# OCFC Yearly Financials - add newest fiscal year column (FY ending 2018-12-31)
# as a new column D, shifting the existing D:K data right to E:L.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OCFC")

# Insert a new blank column before column D; this shifts D:K -> E:L.
$ws.Columns("D:D").Insert()

# The newly inserted column D has no explicit formatting; copy the number
# formats / fonts from column E (the old column D, now shifted one column
# right) so the new column matches the rest of the table exactly.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the newest fiscal year's figures.

# Income Statement
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 276700
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = -3800
$ws.Range("D17").Value = 39600
$ws.Range("D18").Value = 237000
$ws.Range("D20").Value = -151500
$ws.Range("D21").Value = 98100
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 85500
$ws.Range("D24").Value = 15400
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 70100
$ws.Range("D27").Value = 70100
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 1900
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 151500
$ws.Range("D33").Value = 71900
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 71900

# Balance Sheet
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 120800
$ws.Range("D42").Value = 66400
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 111200
$ws.Range("D49").Value = 355400
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 67900
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 7516200
$ws.Range("D57").Value = 0
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 99500
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 6476800
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 305100
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1039400
$ws.Range("D77").Value = 0

# Cash Flow Statement
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 71900
$ws.Range("D83").Value = 12600
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 92600
$ws.Range("D91").Value = -11500
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 48600
$ws.Range("D96").Value = -29600
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -128400
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 12700
